$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated/corrected values for existing rows (columns H = AgTests, I = AgPosit)
$ws.Range("H306").Value = 70735

$ws.Range("H310").Value = 75099
$ws.Range("I310").Value = 3928

$ws.Range("H313").Value = 61399

$ws.Range("H316").Value = 49212

$ws.Range("H320").Value = 75017
$ws.Range("I320").Value = 3643

$ws.Range("H321").Value = 90515
$ws.Range("I321").Value = 2791

$ws.Range("H322").Value = 107315

$ws.Range("H323").Value = 148928
$ws.Range("I323").Value = 2288

$ws.Range("H324").Value = 232599
$ws.Range("I324").Value = 2657

$ws.Range("H325").Value = 727616
$ws.Range("I325").Value = 6003

$ws.Range("H326").Value = 426002
$ws.Range("I326").Value = 3741

$ws.Range("H327").Value = 235696
$ws.Range("I327").Value = 2871

$ws.Range("H328").Value = 178116

$ws.Range("H329").Value = 82171

$ws.Range("H330").Value = 70704
$ws.Range("I330").Value = 1990

$ws.Range("H331").Value = 150106
$ws.Range("I331").Value = 2585

$ws.Range("H332").Value = 420558
$ws.Range("I332").Value = 4084

$ws.Range("H333").Value = 256462

$ws.Range("H334").Value = 201568
$ws.Range("I334").Value = 3377

$ws.Range("H335").Value = 129494
$ws.Range("I335").Value = 2901

$ws.Range("H336").Value = 100221
$ws.Range("I336").Value = 3176

$ws.Range("H337").Value = 101953
$ws.Range("I337").Value = 2879

$ws.Range("H338").Value = 215654
$ws.Range("I338").Value = 3073

$ws.Range("H339").Value = 594288
$ws.Range("I339").Value = 5156

$ws.Range("H340").Value = 331907
$ws.Range("I340").Value = 2978

$ws.Range("H341").Value = 391294
$ws.Range("I341").Value = 4570

# Append new row 342 with the latest daily stats (st 10. 02. 2021)
$ws.Range("A342").Value = 44236
$ws.Range("A342").NumberFormat = $ws.Range("A341").NumberFormat

$ws.Range("B342").Value = 268986
$ws.Range("C342").Value = 251618
$ws.Range("D342").Value = 11866
$ws.Range("E342").Value = 13980
$ws.Range("F342").Value = 3179
$ws.Range("G342").Value = 5502
$ws.Range("H342").Value = 167556
$ws.Range("I342").Value = 3067
